# Applies the "Add files via upload" edit to
# software_update_technical.docx: documents the new
# '/var/www/code_updates' staging directory that the nightly/manual
# update-check scripts now use, and the subsequent copy/cleanup flow.

$d = $word.ActiveDocument

$lsq = [char]0x2018   # left single quote
$rsq = [char]0x2019   # right single quote

# ---------------------------------------------------------------------
# Change 1 (Paragraph 2): "MaxAir executes a script each day ..." --
# mention that discrepant files get downloaded to the staging directory
# before the 'Update' icon is shown.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2).Range
$old1 = "code modules in the Repository. If there are any discrepancies then an " + $lsq + "Update" + $rsq + " icon will be displayed on the toolbar: "
$new1 = "code modules in the Repository. If there are any discrepancies then copies of the updated files will be downloaded to directory " + $lsq + "/var/www/code_updates" + $rsq + " and an " + $lsq + "Update" + $rsq + " icon will be displayed on the toolbar: "
$p2.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2 (Paragraph 9): "Selecting 'Yes' will overwrite ..." --
# describe the overwrite as coming from the staged file/directory, and
# that the staging copies are deleted afterwards.
# ---------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9).Range
$old2 = "software code modules held in the Repository."
$new2 = "file held in the directory " + $lsq + "/var/www/code_updates" + $rsq + ". The copies in " + $lsq + "/var/www/code_updates" + $rsq + " will then be deleted, removing the " + $lsq + "Update" + $rsq + " icon from the toolbar."
$p9.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3 (Paragraph 14): "The current and Repository versions will be
# displayed." -- text itself is unchanged upstream; this just re-applies
# it so the previously split trailing runs normalise into one.
# ---------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14).Range
$old3 = "The current and Repository versions will be displayed."
$new3 = "The current and Repository versions will be displayed."
$p14.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------
# Change 4 (Paragraph 21): "A dialogue will be displayed ..." --
# same staging-directory explanation for the manual 'Check for Updates'
# background task.
# ---------------------------------------------------------------------
$p21 = $d.Paragraphs.Item(21).Range
$old4 = "if there are any discrepancies then the " + $lsq + "Update" + $rsq + " icon will be displayed on the toolbar."
$new4 = "if there are any discrepancies then copies of the updated files will be stored in directory " + $lsq + "/var/www/code_updates" + $rsq + " and the " + $lsq + "Update" + $rsq + " icon will be displayed on the toolbar."
$p21.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

Write-Host "Edit applied."
